$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row (row 1): insert a new "n" header in B1, which pushes the
# existing "Trial N (s)" headers from B:L to C:M... actually the trial
# headers stay in C:L (one column to the right of their old B:K position),
# "Average (s)" stays in M1, and a brand new "Points for Graph" header is
# added in N1.
# ---------------------------------------------------------------------------
$trialHeaders = @("Trial 1 (s)","Trial 2 (s)","Trial 3 (s)","Trial 4 (s)","Trial 5 (s)","Trial 6 (s)","Trial 7 (s)","Trial 8 (s)","Trial 9 (s)","Trial 10 (s)")

$ws.Cells.Item(1, 1).Value2 = "Brute Force"
$ws.Cells.Item(1, 2).Value2 = "n"
for ($i = 0; $i -lt $trialHeaders.Length; $i++) {
    $col = 3 + $i  # C=3 .. L=12
    $ws.Cells.Item(1, $col).Value2 = $trialHeaders[$i]
}
$ws.Cells.Item(1, 13).Value2 = "Average (s)"
$ws.Cells.Item(1, 13).Font.Bold = $true
$ws.Cells.Item(1, 14).Value2 = "Points for Graph"
$ws.Cells.Item(1, 14).Font.Bold = $true

# ---------------------------------------------------------------------------
# "n" blocks for Brute Force (rows 2-6) and Naive DNC (rows 8-12): drop the
# old "n = ..." text labels from column A entirely and write plain numeric
# n values into column B instead.
# ---------------------------------------------------------------------------
$nValues = @(10, 100, 1000, 10000, 100000)

for ($i = 0; $i -lt $nValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Clear()
    $ws.Cells.Item($row, 2).Value2 = $nValues[$i]
}

$ws.Cells.Item(7, 1).Value2 = "Naive DNC"
$ws.Cells.Item(7, 2).Font.Bold = $true

for ($i = 0; $i -lt $nValues.Length; $i++) {
    $row = 8 + $i
    $ws.Cells.Item($row, 1).Clear()
    $ws.Cells.Item($row, 2).Value2 = $nValues[$i]
}

$ws.Cells.Item(13, 1).Value2 = "Enhanced DNC"
$ws.Cells.Item(13, 2).Font.Bold = $true

# ---------------------------------------------------------------------------
# Data rows 14-18: shift trial values from B:K to C:L, put the numeric n in
# B, fix the average formula to reference C:L, and add the new
# "Points for Graph" formula in column N.
# ---------------------------------------------------------------------------
$dataRows = @(14, 15, 16, 17, 18)
$nForRow = @{ 14 = 10; 15 = 100; 16 = 1000; 17 = 10000; 18 = 100000 }

foreach ($row in $dataRows) {
    # Capture old trial values (B:K) before overwriting anything, using
    # Value2 (raw numeric) rather than Value/Formula (which stringify
    # scientific notation and trip up automatic number-format detection).
    $old = @()
    for ($c = 2; $c -le 11; $c++) {
        $old += , ($ws.Cells.Item($row, $c).Value2)
    }

    # Clear the old "n = ..." label that lived in column A.
    $ws.Cells.Item($row, 1).Clear()

    # Shift the captured trial values right into C:L.
    for ($i = 0; $i -lt $old.Length; $i++) {
        $ws.Cells.Item($row, 3 + $i).Value2 = $old[$i]
    }

    # Plain numeric n value now lives in column B.
    $ws.Cells.Item($row, 2).Value2 = $nForRow[$row]

    # Average formula now spans C:L instead of B:K.
    $ws.Cells.Item($row, 13).Formula = "=AVERAGE(C$row" + ":L$row)"

    # New "Points for Graph" column: textual "(n, average)" pair.
    $ws.Cells.Item($row, 14).Formula = '="(" & B' + $row + ' & ", " & M' + $row + ' & ")"'
}

# ---------------------------------------------------------------------------
# Column widths: column B gets the (nearly-default) 11.53 width, column N is
# widened to fit the new "(n, avg)" strings.
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 11.53
$ws.Columns.Item(14).ColumnWidth = 17.95

# ---------------------------------------------------------------------------
# An extra, otherwise-empty row 24 (column M) shows up in the saved file —
# touch it so the used range extends down to row 24.
# ---------------------------------------------------------------------------
$ws.Cells.Item(24, 13).NumberFormat = "General"

# ---------------------------------------------------------------------------
# Selection matches the post-edit state captured in the diff.
# ---------------------------------------------------------------------------
$ws.Range("N14:N18").Select()
